$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 141 - this shifts rows 141:166 down to 142:167
$ws.Rows.Item(141).Insert()

# Fill in the newly inserted row 141 with the new data record
$ws.Cells.Item(141, 1).Value = 10
$ws.Cells.Item(141, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(141, 3).Value = "La Araucanía"
$ws.Cells.Item(141, 4).Value = 44522
$ws.Cells.Item(141, 5).Value = 9
$ws.Cells.Item(141, 6).Value = "Fruta"
$ws.Cells.Item(141, 7).Value = 100102
$ws.Cells.Item(141, 8).Value = "Cítricos"
$ws.Cells.Item(141, 9).Value = 100102006
$ws.Cells.Item(141, 10).Value = "Pomelo"
$ws.Cells.Item(141, 11).Value = "Start Ruby"
$ws.Cells.Item(141, 12).Value = "Segunda"
$ws.Cells.Item(141, 13).Value = 40
$ws.Cells.Item(141, 14).Value = 8000
$ws.Cells.Item(141, 15).Value = 8000
$ws.Cells.Item(141, 16).Value = 8000
$ws.Cells.Item(141, 17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item(141, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(141, 19).Value = 533
$ws.Cells.Item(141, 20).Value = 15
